{"js": "// Apply the billing-document edit described by the diff:\n//  - Update the first invoice block's retailer/address/phone/GSTIN/beat/\n//    store-id/time/amount-in-words/summary-line fields for the new\n//    \"ARVIND AGENCY\" retailer.\n//  - Remove the entire second invoice block (the \"DEVAKI ENTERPRISES\" /\n//    \"JAYALAXMI STORE...\" copy that followed the page break run), which\n//    was dropped from the document.\n\nconst body = context.document.body;\n\n// Exact-text replacements (each `oldText` is unique in the document, so a\n// single search-and-replace on the whole paragraph text is unambiguous).\nconst replacements = [\n  [\n    \"  ARIYAMANGALAM                                      Retailer Name   :  ABU MALIGAI-D                      \",\n    \"  ARIYAMANGALAM                                      Retailer Name   :  ARVIND AGENCY                      \",\n  ],\n  [\n    \"  TRICHY-620010                                      Address         :  Marsingpet Road                    \",\n    \"  TRICHY-620010                                      Address         :  2/3, Nadu ValayalKara Street       \",\n  ],\n  [\n    \"  PHONE NO         :9944951444                                          BEEMANAGAR                         \",\n    \"  PHONE NO         :9944951444                                          Tharanallur, Tiruchirappalli,      \",\n  ],\n  [\n    \"  GSTIN No         :33AAPFD1365C1ZR                                     Trichy                             \",\n    \"  GSTIN No         :33AAPFD1365C1ZR                                     India                              \",\n  ],\n  [\n    \"  RS PAN No        :AAPFD1365C                       Phone No        :   8072948180                        \",\n    \"  RS PAN No        :AAPFD1365C                       Phone No        :   6369137944                        \",\n  ],\n  [\n    \"  Beat Name        :D-BEEMANAGAR                     GSTIN NO        :                                     \",\n    \"  Beat Name        :D-WHOLESALE                      GSTIN NO        :   33BQTPM0121J1ZW                   \",\n  ],\n  [\n    \"  HUL STORE ID     :HUL-41A392D-P18479               Time of Billing :   01/04/2023 10:35:49               \",\n    \"  HUL STORE ID     :HUL-41A392D-P25102               Time of Billing :   01/04/2023 10:19:36               \",\n  ],\n  [\n    \"  Fifteen Thousand Five Hundred Fifty-Three Rupees Only                      \",\n    \"  Fifty-Five Thousand Seven Rupees Only                                      \",\n  ],\n  [\n    \"  ABC00003     ABU MALIGAI-D    Amt : 15553.00\",\n    \"  ABC00001     ARVIND AGENCY    Amt : 55007.00\",\n  ],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly one match for ${JSON.stringify(oldText)}, found ${results.items.length}`\n    );\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// Remove the whole second invoice block: the 15 paragraphs starting at the\n// second \"DEVAKI ENTERPRISES\" paragraph through the trailing page-break\n// paragraph (the last paragraph of the body). Deleting paragraph-by-paragraph\n// (back to front) avoids leaving a stray empty paragraph behind, which a\n// single combined-range delete would do since the body must always end with\n// a paragraph mark.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nlet blockStart = -1;\nlet devakiSeen = 0;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.trim() === \"DEVAKI ENTERPRISES\") {\n    devakiSeen++;\n    if (devakiSeen === 2) {\n      blockStart = i;\n      break;\n    }\n  }\n}\n\nif (blockStart === -1) {\n  throw new Error(\"Could not find the second 'DEVAKI ENTERPRISES' paragraph to remove.\");\n}\n\nfor (let i = items.length - 1; i >= blockStart; i--) {\n  items[i].delete();\n}\nawait context.sync();\n", "ps1": "# Apply the billing-document edit described by the diff:\n#  - Update the first invoice block's retailer/address/phone/GSTIN/beat/\n#    store-id/time/amount-in-words/summary-line fields for the new\n#    \"ARVIND AGENCY\" retailer.\n#  - Remove the entire second invoice block (the \"DEVAKI ENTERPRISES\" /\n#    \"JAYALAXMI STORE...\" copy that followed the page break run), which\n#    was dropped from the document.\n\n$d = $word.ActiveDocument\n\nfunction Replace-ExactText($oldText, $newText) {\n    $range = $d.Content\n    $find = $range.Find\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 1)\n    if (-not $found) {\n        throw \"Could not find text to replace: $oldText\"\n    }\n}\n\nReplace-ExactText \"  ARIYAMANGALAM                                      Retailer Name   :  ABU MALIGAI-D                      \" \"  ARIYAMANGALAM                                      Retailer Name   :  ARVIND AGENCY                      \"\nReplace-ExactText \"  TRICHY-620010                                      Address         :  Marsingpet Road                    \" \"  TRICHY-620010                                      Address         :  2/3, Nadu ValayalKara Street       \"\nReplace-ExactText \"  PHONE NO         :9944951444                                          BEEMANAGAR                         \" \"  PHONE NO         :9944951444                                          Tharanallur, Tiruchirappalli,      \"\nReplace-ExactText \"  GSTIN No         :33AAPFD1365C1ZR                                     Trichy                             \" \"  GSTIN No         :33AAPFD1365C1ZR                                     India                              \"\nReplace-ExactText \"  RS PAN No        :AAPFD1365C                       Phone No        :   8072948180                        \" \"  RS PAN No        :AAPFD1365C                       Phone No        :   6369137944                        \"\nReplace-ExactText \"  Beat Name        :D-BEEMANAGAR                     GSTIN NO        :                                     \" \"  Beat Name        :D-WHOLESALE                      GSTIN NO        :   33BQTPM0121J1ZW                   \"\nReplace-ExactText \"  HUL STORE ID     :HUL-41A392D-P18479               Time of Billing :   01/04/2023 10:35:49               \" \"  HUL STORE ID     :HUL-41A392D-P25102               Time of Billing :   01/04/2023 10:19:36               \"\nReplace-ExactText \"  Fifteen Thousand Five Hundred Fifty-Three Rupees Only                      \" \"  Fifty-Five Thousand Seven Rupees Only                                      \"\nReplace-ExactText \"  ABC00003     ABU MALIGAI-D    Amt : 15553.00\" \"  ABC00001     ARVIND AGENCY    Amt : 55007.00\"\n\n# Find the second \"DEVAKI ENTERPRISES\" paragraph (start of the duplicated\n# second invoice block) and remove everything from there to the end of the\n# document (including the trailing page-break paragraph).\n$devakiCount = 0\n$blockStartParagraph = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $paraText = $d.Paragraphs.Item($i).Range.Text.Trim()\n    if ($paraText -eq \"DEVAKI ENTERPRISES\") {\n        $devakiCount++\n        if ($devakiCount -eq 2) {\n            $blockStartParagraph = $d.Paragraphs.Item($i)\n            break\n        }\n    }\n}\n\nif ($blockStartParagraph -eq $null) {\n    throw \"Could not find the second 'DEVAKI ENTERPRISES' paragraph to remove.\"\n}\n\n$deleteRange = $d.Range($blockStartParagraph.Range.Start, $d.Content.End)\n$deleteRange.Delete()\n"}
